$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "25.821.07"
Set-TextValue "E2" "  -0.06%  "
Set-TextValue "D3" "1.636.17"
Set-TextValue "E3" "  -0.04%  "
Set-TextValue "E4" "  -0.14%  "
Set-TextValue "D5" "215.43"
Set-TextValue "E5" "  -0.19%  "
Set-TextValue "D6" "0.504"
Set-TextValue "E6" "  -0.39%  "
Set-TextValue "E7" "  -0.08%  "
Set-TextValue "E8" "  +0.22%  "
Set-TextValue "E9" "  +0.07%  "
Set-TextValue "D10" "19.88"
Set-TextValue "E10" "  +1.64%  "
Set-TextValue "E11" "  +0.77%  "
Set-TextValue "E12" "  -0.77%  "
Set-TextValue "D13" "1.641.89"
Set-TextValue "E13" "  +0.27%  "
Set-TextValue "D14" "1.860.80"
Set-TextValue "E15" "  -0.65%  "
Set-TextValue "D16" "0.0₃0772"
Set-TextValue "E16" "  +1.72%  "
Set-TextValue "D17" "63.15"
Set-TextValue "E17" "  -0.08%  "
Set-TextValue "D18" "25.825.76"
Set-TextValue "E19" "  -0.12%  "
Set-TextValue "D20" "4.44"
Set-TextValue "E20" "  +2.63%  "
Set-TextValue "D21" "194.46"
Set-TextValue "E21" "  +0.06%  "
Set-TextValue "E22" "  +0.86%  "
Set-TextValue "E23" "  +1.22%  "
Set-TextValue "E24" "  -0.04%  "
Set-TextValue "D25" "1.78"
Set-TextValue "E25" "  -0.79%  "
Set-TextValue "D26" "139.00"
Set-TextValue "E26" "  -0.80%  "
Set-TextValue "D27" "0.122"
Set-TextValue "E27" "  -4.57%  "
Set-TextValue "D28" "6.85"
Set-TextValue "E28" "  +1.07%  "
Set-TextValue "D29" "15.57"
Set-TextValue "E29" "  +0.81%  "
Set-TextValue "E30" "  +0.29%  "
Set-TextValue "D31" "0.0499"
Set-TextValue "E31" "  +2.28%  "
Set-TextValue "E32" "  +1.13%  "
Set-TextValue "E33" "  +1.39%  "
Set-TextValue "E34" "  +2.30%  "
Set-TextValue "E35" "  +0.82%  "
Set-TextValue "D36" "0.901"
Set-TextValue "E36" "  -0.11%  "
Set-TextValue "E37" "  -0.25%  "
Set-TextValue "D38" "0.552"
Set-TextValue "E38" "  +0.52%  "
Set-TextValue "D39" "1.109.09"
Set-TextValue "E39" "  -1.82%  "
Set-TextValue "E40" "  +0.30%  "
Set-TextValue "E41" "  +0.11%  "
Set-TextValue "E42" "  +0.70%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "0.804"
Set-TextValue "E43" "  +0.39%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D44" "99.31"
Set-TextValue "E44" "  +1.92%  "
Set-TextValue "D45" "0.0₆0110"
Set-TextValue "E45" "  -3.36%  "
Set-TextValue "D46" "55.59"
Set-TextValue "E46" "  +0.20%  "
Set-TextValue "D47" "2.51"
Set-TextValue "E47" "  +12.08%  "
Set-TextValue "D48" "0.418"
Set-TextValue "E48" "  -5.88%  "
Set-TextValue "D49" "7.69"
Set-TextValue "E49" "  -0.14%  "
Set-TextValue "E50" "  -0.33%  "
Set-TextValue "D51" "1.00"
Set-TextValue "E51" "  -0.09%  "
